$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: In the opening "Law and storytelling..." paragraph, change
#   "each taking the floor and doing their best"
# to
#   'each "take the floor" and do their best'
#
# Note: We locate the text with Find.Execute (no replacement argument) and
# then assign Range.Text directly, rather than passing a replacement string
# to Find.Execute. Using Find.Execute's replacement argument triggers the
# "smart quotes" AutoFormat/AutoCorrect behavior which turns straight quotes
# into curly quotes; assigning Range.Text inserts the literal straight
# quotes we want, matching the source diff exactly.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("each taking the floor and doing their best")
if ($found1) {
    $rng1.Text = 'each "take the floor" and do their best'
} else {
    Write-Output "WARNING: change 1 target text not found"
}

# ---------------------------------------------------------------------------
# Change 2: In the Scylla/Charybdis paragraph, add a comma after "strait":
#   "on the Italian side of the strait and Charybdis"
# to
#   "on the Italian side of the strait, and Charybdis"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("on the Italian side of the strait and Charybdis")
if ($found2) {
    $rng2.Text = "on the Italian side of the strait, and Charybdis"
} else {
    Write-Output "WARNING: change 2 target text not found"
}
